# Actualizacion automatica del tracker
# Appends the latest tracked results to the bottom of the results table.
#
# NOTE: this runtime's PowerShell does not bind named parameters
# (`-Row 76`) correctly inside function bodies - the parameter reads back
# empty. Positional parameters work fine, so every helper call below uses
# positional args.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TrackerRow {
    param([int]$Row, [double]$EventId, [string]$Fecha, [string]$JugadorA, [string]$JugadorB, [string]$Pronostico, [double]$Cuota)

    $ws.Range("A$Row").Value = $EventId

    # Prefix with a single quote so the YYYY-MM-DD text isn't coerced into a
    # serial date value - the tracker stores "fecha" as plain text, and then
    # strip the resulting quote-prefix formatting so no style is left on the
    # cell (matches the rest of the column).
    $ws.Range("B$Row").Value = "'" + $Fecha
    $ws.Range("B$Row").ClearFormats()

    $ws.Range("C$Row").Value = $JugadorA
    $ws.Range("D$Row").Value = $JugadorB
    $ws.Range("E$Row").Value = $Pronostico
    $ws.Range("F$Row").Value = $Cuota

    # resultado / profit are still pending for these freshly-added matches.
    # The tracker represents "pending" as an explicit empty TEXT cell, not a
    # blank cell - plain "" would clear the cell entirely, so force a
    # text-typed empty value via the same quote-prefix trick, then clear the
    # formatting it leaves behind.
    $ws.Range("G$Row").Value = "'"
    $ws.Range("G$Row").ClearFormats()
    $ws.Range("H$Row").Value = "'"
    $ws.Range("H$Row").ClearFormats()
}

Set-TrackerRow 76 14581060 "2025-09-02" "Pol Martin Tiffon" "Henrique Rocha" "Gana Pol Martin Tiffon" 3.25
Set-TrackerRow 77 14592779 "2025-09-02" "Robert Strombachs" "Frederico Ferreira Silva" "Gana Robert Strombachs" 3.25
Set-TrackerRow 78 14592266 "2025-09-02" "Kokoro Isomura" "Renta Tokuda" "Gana Kokoro Isomura" 2.75

Write-Host "Tracker actualizado: filas 76-78 agregadas."
